$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "From"
$ws.Range("B1").Value = "To"
$ws.Range("D1").Value = "Dragging Comment"
$ws.Range("C1").Value = "Date"

$ws.Range("A2").Value = "FzConfig"
$ws.Range("B2").Value = "ActConfig"
$ws.Range("D2").Value = "tested running; appeared fine"

$ws.Range("A3").Select()
